# Scheduled runner update: refresh market-board derived price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 64013.375
$ws.Range("J17").Value = 72993.86
$ws.Range("L17").Value = 218981.58
$ws.Range("N17").Value = -219317.58
# Row 70
$ws.Range("H70").Value = 1356713.8
$ws.Range("J70").Value = 2000.4286
$ws.Range("L70").Value = 6001.2858
$ws.Range("N70").Value = -6541.2858
# Row 73
$ws.Range("H73").Value = 1356713.8
$ws.Range("J73").Value = 2000.4286
$ws.Range("L73").Value = 6001.2858
$ws.Range("N73").Value = -7873.2858
# Row 88
$ws.Range("H88").Value = 1970.125
$ws.Range("I88").Value = 1003
$ws.Range("J88").Value = 2108.2856
$ws.Range("K88").Value = 1003
$ws.Range("L88").Value = 2108.2856
$ws.Range("M88").Value = -597
$ws.Range("N88").Value = -2920.2856
# Row 91
$ws.Range("H91").Value = 1970.125
$ws.Range("I91").Value = 1003
$ws.Range("J91").Value = 2108.2856
$ws.Range("K91").Value = 1003
$ws.Range("L91").Value = 2108.2856
$ws.Range("M91").Value = 401
$ws.Range("N91").Value = -4916.2856
# Row 101
$ws.Range("H101").Value = 618.8
$ws.Range("J101").Value = 408.76923
$ws.Range("L101").Value = 1226.30769
$ws.Range("N101").Value = -4470.30769
# Row 131
$ws.Range("H131").Value = 1626208.4
$ws.Range("I131").Value = 1948.909
$ws.Range("J131").Value = 7581826.5
$ws.Range("K131").Value = 5846.727000000001
$ws.Range("L131").Value = 22745479.5
$ws.Range("M131").Value = -806.7270000000008
$ws.Range("N131").Value = -22755559.5
# Row 132
$ws.Range("H132").Value = 2613.1155
$ws.Range("I132").Value = 2715.4783
$ws.Range("J132").Value = 1828.3334
$ws.Range("K132").Value = 8146.4349
$ws.Range("L132").Value = 5485.0002
$ws.Range("M132").Value = -5616.4349
$ws.Range("N132").Value = -10545.0002
# Row 135
$ws.Range("H135").Value = 2862.2122
$ws.Range("I135").Value = 418.2
$ws.Range("K135").Value = 3763.8
$ws.Range("M135").Value = -1228.8
# Row 138
$ws.Range("H138").Value = 3443.8225
$ws.Range("I138").Value = 2295.7917
$ws.Range("J138").Value = 4168.8945
$ws.Range("K138").Value = 6887.375100000001
$ws.Range("L138").Value = 12506.6835
$ws.Range("M138").Value = -1747.375100000001
$ws.Range("N138").Value = -22786.6835

$ws = $wb.Worksheets.Item("ARM")
# Row 60
$ws.Range("H60").Value = 89473.336
$ws.Range("I60").Value = 89473.336
$ws.Range("K60").Value = 89473.336
$ws.Range("M60").Value = -88740.336
# Row 61
$ws.Range("H61").Value = 3129054.5
$ws.Range("I61").Value = 3575386.5
$ws.Range("J61").Value = 4730.5
$ws.Range("K61").Value = 3575386.5
$ws.Range("L61").Value = 4730.5
$ws.Range("M61").Value = -3575174.5
$ws.Range("N61").Value = -5154.5
# Row 74
$ws.Range("H74").Value = 2296.56
$ws.Range("I74").Value = 1972.5264
$ws.Range("K74").Value = 1972.5264
$ws.Range("M74").Value = -1098.5264
# Row 77
$ws.Range("H77").Value = 2296.56
$ws.Range("I77").Value = 1972.5264
$ws.Range("K77").Value = 9862.632
$ws.Range("M77").Value = -5494.632
# Row 136
$ws.Range("H136").Value = 3129054.5
$ws.Range("I136").Value = 3575386.5
$ws.Range("J136").Value = 4730.5
$ws.Range("K136").Value = 10726159.5
$ws.Range("L136").Value = 14191.5
$ws.Range("M136").Value = -10723609.5
$ws.Range("N136").Value = -19291.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 19611576
$ws.Range("I31").Value = 37039710
$ws.Range("J31").Value = 4922.7085
$ws.Range("K31").Value = 37039710
$ws.Range("L31").Value = 4922.7085
$ws.Range("M31").Value = -37039415
$ws.Range("N31").Value = -5512.7085
# Row 34
$ws.Range("H34").Value = 19611576
$ws.Range("I34").Value = 37039710
$ws.Range("J34").Value = 4922.7085
$ws.Range("K34").Value = 37039710
$ws.Range("L34").Value = 4922.7085
$ws.Range("M34").Value = -37039508
$ws.Range("N34").Value = -5326.7085
# Row 132
$ws.Range("H132").Value = 1753.4359
$ws.Range("I132").Value = 1740
$ws.Range("K132").Value = 5220
$ws.Range("M132").Value = -2690
# Row 134
$ws.Range("H134").Value = 2208.682
$ws.Range("I134").Value = 1951.7858
$ws.Range("K134").Value = 5855.357400000001
$ws.Range("M134").Value = -3320.357400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 649.3333
$ws.Range("I5").Value = 582.7143
$ws.Range("J5").Value = 782.5714
$ws.Range("K5").Value = 1748.1429
$ws.Range("L5").Value = 2347.7142
$ws.Range("M5").Value = -1636.1429
$ws.Range("N5").Value = -2571.7142
# Row 135
$ws.Range("H135").Value = 649.3333
$ws.Range("I135").Value = 582.7143
$ws.Range("J135").Value = 782.5714
$ws.Range("K135").Value = 5244.428699999999
$ws.Range("L135").Value = 7043.1426
$ws.Range("M135").Value = -2709.428699999999
$ws.Range("N135").Value = -12113.1426

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2447.077
$ws.Range("I80").Value = 1931.5
$ws.Range("K80").Value = 1931.5
$ws.Range("M80").Value = -933.5
# Row 83
$ws.Range("H83").Value = 2447.077
$ws.Range("I83").Value = 1931.5
$ws.Range("K83").Value = 9657.5
$ws.Range("M83").Value = -4665.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2085338.2
$ws.Range("I68").Value = 2977754.8
$ws.Range("J68").Value = 3033.3333
$ws.Range("K68").Value = 2977754.8
$ws.Range("L68").Value = 3033.3333
$ws.Range("M68").Value = -2977005.8
$ws.Range("N68").Value = -4531.3333
# Row 71
$ws.Range("H71").Value = 2085338.2
$ws.Range("I71").Value = 2977754.8
$ws.Range("J71").Value = 3033.3333
$ws.Range("K71").Value = 14888774
$ws.Range("L71").Value = 15166.6665
$ws.Range("M71").Value = -14885030
$ws.Range("N71").Value = -22654.6665
# Row 82
$ws.Range("H82").Value = 6727.857
$ws.Range("I82").Value = 3688
$ws.Range("J82").Value = 12199.6
$ws.Range("K82").Value = 3688
$ws.Range("L82").Value = 12199.6
$ws.Range("M82").Value = -3327
$ws.Range("N82").Value = -12921.6
# Row 85
$ws.Range("H85").Value = 6727.857
$ws.Range("I85").Value = 3688
$ws.Range("J85").Value = 12199.6
$ws.Range("K85").Value = 3688
$ws.Range("L85").Value = 12199.6
$ws.Range("M85").Value = -2440
$ws.Range("N85").Value = -14695.6

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
# Row 43
$ws.Range("H43").Value = 89999
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
# Row 62
$ws.Range("H62").Value = 1225107.2
$ws.Range("J62").Value = 2287327
$ws.Range("L62").Value = 2287327
$ws.Range("N62").Value = -2288575
# Row 65
$ws.Range("H65").Value = 1225107.2
$ws.Range("J65").Value = 2287327
$ws.Range("L65").Value = 11436635
$ws.Range("N65").Value = -11442875
# Row 107
$ws.Range("H107").Value = 3896.282
$ws.Range("I107").Value = 2339.0908
$ws.Range("J107").Value = 5911.4707
$ws.Range("K107").Value = 7017.2724
$ws.Range("L107").Value = 17734.4121
$ws.Range("M107").Value = -5097.2724
$ws.Range("N107").Value = -21574.4121
# Row 113
$ws.Range("H113").Value = 540.8889
$ws.Range("I113").Value = 478.25
$ws.Range("J113").Value = 719.8570999999999
$ws.Range("K113").Value = 1434.75
$ws.Range("L113").Value = 2159.5713
$ws.Range("M113").Value = 735.25
$ws.Range("N113").Value = -6499.5713
# Row 132
$ws.Range("H132").Value = 1076.2222
$ws.Range("I132").Value = 958.6667
$ws.Range("J132").Value = 1593.4667
$ws.Range("K132").Value = 2876.0001
$ws.Range("L132").Value = 4780.4001
$ws.Range("M132").Value = -346.0001000000002
$ws.Range("N132").Value = -9840.400099999999
# Row 136
$ws.Range("H136").Value = 2534.9424
$ws.Range("I136").Value = 2292.4773
$ws.Range("J136").Value = 3868.5
$ws.Range("K136").Value = 6877.4319
$ws.Range("L136").Value = 11605.5
$ws.Range("M136").Value = -4327.4319
$ws.Range("N136").Value = -16705.5

Write-Host "Applied scheduled price/profit refresh to ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets."
